# Fruta / hortaliza, semanal
# The data rows (2-20) get reshuffled: each destination row receives the
# D, J, K, L, M, O, P values that used to belong to a different source row.
# All other columns (A,B,C,E,F,G,H,I,N,Q,R) are identical across every row,
# so they do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row number -> source row number (values to copy from).
$mapping = @{}
$mapping[2]  = 14
$mapping[3]  = 10
$mapping[4]  = 9
$mapping[5]  = 8
$mapping[6]  = 3
$mapping[7]  = 17
$mapping[8]  = 16
$mapping[9]  = 13
$mapping[10] = 2
$mapping[11] = 15
$mapping[12] = 18
$mapping[13] = 5
$mapping[14] = 19
$mapping[15] = 20
$mapping[16] = 7
$mapping[17] = 11
$mapping[18] = 4
$mapping[19] = 12
$mapping[20] = 6

# Snapshot the current (pre-edit) values for the columns that move, so that
# writing into the destination rows never clobbers a value still needed as
# a source for a later write (the mapping is a permutation).
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

for ($destRow = 2; $destRow -le 20; $destRow++) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value2  = $src.D   # D: Fecha
    $ws.Cells.Item($destRow, 10).Value2 = $src.J   # J: Volumen
    $ws.Cells.Item($destRow, 11).Value2 = $src.K   # K: Precio minimo
    $ws.Cells.Item($destRow, 12).Value2 = $src.L   # L: Precio maximo
    $ws.Cells.Item($destRow, 13).Value2 = $src.M   # M: Precio promedio ponderado
    $ws.Cells.Item($destRow, 15).Value2 = $src.O   # O: Origen
    $ws.Cells.Item($destRow, 16).Value2 = $src.P   # P: Precio $/Kg
}
